# Update automatico via Actualizar 03-04-2021 23-13-01
#
# This refreshes the "Disponibilidad" log sheet: the newest timestamp block
# (rows 2-15) gets a newer "Fecha" stamp, the block that used to be newest
# (rows 16-29) shifts to become the "Ultimo" historical block with the prior
# rows' stamp, and the oldest block (rows 30-43) is overwritten with the
# values the middle block used to hold (the very oldest snapshot falls off).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Most-recent block (rows 2-15): refresh the "Fecha" timestamp ---------
$ws.Range("D2:D15").Value = 44259.96680961041

# Row 4 previously carried a one-off "Ultimo" stamp in E4 that is no longer
# part of this snapshot - remove the cell entirely.
$ws.Range("E4").Clear()

# --- Second block (rows 16-29): becomes the new "previous" snapshot -------
$ws.Range("D16:D29").Value = 44259.94453363426

# Row 18 (Shiny) flips back to "Disponible" and records an "Ultimo" stamp.
$ws.Range("C18").Value = "Disponible"
$ws.Range("E18").Value = 44259.94448752315
$ws.Range("E18").NumberFormat = $ws.Range("D18").NumberFormat

# --- Oldest block (rows 30-43): overwritten with the prior block's stamp --
$ws.Range("D30:D43").Value = 44250.35508177083

# Row 32 (Shiny) flips to "No Disponible" for this older snapshot.
$ws.Range("C32").Value = "No Disponible"
